# Interdiff between v17 and v18
#
# Replace the "Up Arrow 5" autoshape (which pointed from the
# "ab0:AddressBook" table up to the "currentStatePointer = 0" label)
# with a red "Straight Arrow Connector" running the same span.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the old up-arrow autoshape.
$s.Shapes.Item("Up Arrow 5").Delete()

# EMU -> point helper (1 pt = 12700 EMU) so the new connector lands on
# the exact coordinates used by the target diagram.
$emuPerPt = 12700.0
$left   = 1577130 / $emuPerPt
$top    = 3347207 / $emuPerPt
$bottom = (3347207 + 706873) / $emuPerPt

# msoConnectorStraight = 1. Build it top-to-bottom first (so the height
# comes out exact), then flip it vertically so the arrowhead ends up on
# the top (matching flipV="1" in the target XML).
$conn = $s.Shapes.AddConnector(1, $left, $top, $left, $bottom)
$conn.Name = "Straight Arrow Connector 2"

# Zero out any floating point drift on the width (the connector is
# perfectly vertical).
$conn.Width = 0

# Flip vertically (msoFlipVertical = 1).
$conn.Flip(1)

# Red line (RRGGBB 0xC00000) with a triangle arrowhead at the tail end.
$conn.Line.ForeColor.RGB = 0xC0
$conn.Line.EndArrowheadStyle = 2
